$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixer2-BOM")

# Fill in "Have" (C) and "Bought" (D) quantities for several rows.
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 8
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 10
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 20
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 2
$ws.Range("C18").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("C21").Value = 0

# Move the active selection on the sheet.
$ws.Range("D6").Select()
